# The workbook has a sheet named "2024" that tracks monthly transaction
# details/dates. A new September transaction ("latest transaction pan" at
# 2024-09-12 12:22:12) was recorded, which pushes all the existing rows
# (from row 37, the most-recent-so-far September entry, down through the
# end of the sheet at row 140) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row above row 37; this shifts rows 37-140 down to 38-141
# and updates the sheet dimension automatically.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new transaction entry.
$ws.Range("R37").Value = "latest transaction pan"
$ws.Range("S37").Value = "2024-09-12 12:22:12"
